$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8652793169021606
$ws.Range("B1").Value = 1.802945256233215
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.530941963195801
$ws.Range("E1").Value = 1.002193331718445
